# Add results and unfolding with 100 keV threshold
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Beta" row (row 2) - recomputed reconstructed dose columns (F..N)
$ws.Range("F2").Value2 = 12.06299460600843
$ws.Range("G2").Value2 = 10.99641073747695
$ws.Range("H2").Value2 = 13.12313884420613
$ws.Range("I2").Value2 = 1.892620994339465
$ws.Range("J2").Value2 = 1.670871214393642
$ws.Range("K2").Value2 = 2.098504397666163
$ws.Range("L2").Value2 = 0.1513826561181574
$ws.Range("M2").Value2 = 0.1352322530838551
$ws.Range("N2").Value2 = 0.166784074039164

# Update existing "Gamma" row (row 3) - recomputed reconstructed dose columns (F..N)
$ws.Range("F3").Value2 = 0.002979199992955184
$ws.Range("G3").Value2 = 0.001225734052259997
$ws.Range("H3").Value2 = 0.005291759204239123
$ws.Range("I3").Value2 = 0.002731003867812408
$ws.Range("J3").Value2 = 0.001127675264114317
$ws.Range("K3").Value2 = 0.004844346623776378
$ws.Range("L3").Value2 = 0.003056091510215558
$ws.Range("M3").Value2 = 0.001295684125663276
$ws.Range("N3").Value2 = 0.005373240507213315

# Add new "Beta + Gamma" row (row 4), copying the styled index-cell format from row 2
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("A4").Value2 = 2
$ws.Range("B4").Value2 = "Beta + Gamma"
$ws.Range("C4").Value2 = 12.00687180793019
$ws.Range("D4").Value2 = 1.974477778970852
$ws.Range("E4").Value2 = 0.1537386519519979
$ws.Range("F4").Value2 = 12.06597380600138
$ws.Range("G4").Value2 = 10.99763647152921
$ws.Range("H4").Value2 = 13.12843060341037
$ws.Range("I4").Value2 = 1.895351998207277
$ws.Range("J4").Value2 = 1.671998889657756
$ws.Range("K4").Value2 = 2.10334874428994
$ws.Range("L4").Value2 = 0.154438747628373
$ws.Range("M4").Value2 = 0.1365279372095184
$ws.Range("N4").Value2 = 0.1721573145463773
